# Daily attendance processing - 2025-11-14 10:50:21
# Normalises the "Recorded By" column (G): rotates the comma-separated list
# of recorders so the last-listed recorder for each session becomes first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = @($text -split ",\s*")

    if ($parts.Count -gt 1) {
        $rotated = @($parts[-1]) + $parts[0..($parts.Count - 2)]
        $cell.Value = [string]::Join(", ", $rotated)
    }
}
